$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Formulas / values -------------------------------------------------

$ws.Range("I1").Formula = "=0.39*(89.64*0.4+0.31*0.35)"

$ws.Range("N2").Value = 49.55
$ws.Range("N3").Value = 706.26

$ws.Range("I4").Formula = "= 35.52 + 14.03"
$ws.Range("K4").Value = 3341.8
$ws.Range("L4").Formula = "= 49.55+ 706.26+240.13+5.52+2354.21"
$ws.Range("N4").Value = 240.13

$ws.Range("K5").Value = 0.048
$ws.Range("L5").Formula = "=89.95/1854.72"
$ws.Range("N5").Value = 5.52

$ws.Range("I6").Value = 89.95
$ws.Range("N6").Value = 2354.21

$ws.Range("K7").Value = 4491.3999999999996
$ws.Range("L7").Formula = "= 3758.4+751.68"
$ws.Range("M7").Value = 3341.8
$ws.Range("N7").Formula = "=SUM(N2:N6)"

$ws.Range("K8").Value = 401.02
$ws.Range("L8").Formula = "=12*3355.67/100"
$ws.Range("O8").Formula = "=L8/L7"

$ws.Range("I9").Value = 619.53
$ws.Range("K9").Value = 3742.82
$ws.Range("L9").Formula = "= 3355.67 + 402.68"
$ws.Range("O9").Formula = "=14.03/4510.1"

$ws.Range("I10").Value = 706.26
$ws.Range("K10").Value = 748.56
$ws.Range("L10").Formula = "=20/100*3758.4"
$ws.Range("O10").Formula = "=L4/L7"

$ws.Range("I11").Value = 240.13

$ws.Range("I12").Formula = "=380*I9/100"

# --- Formatting ----------------------------------------------------------
# Group A: yellow fill + Times New Roman 14pt, no border
$groupA = $ws.Range("I1,I4,L4,L7,L8,L9")
$groupA.Interior.Color = 65535
$groupA.Font.Name = "Times New Roman"
$groupA.Font.Size = 14

# Group B: yellow fill, default font (Calibri 11), no border
$groupB = $ws.Range("I6,N7,I10,I11,I12")
$groupB.Interior.Color = 65535

# Group C: yellow fill + Times New Roman 12pt, centered/wrapped
$groupC = $ws.Range("I9")
$groupC.Interior.Color = 65535
$groupC.Font.Name = "Times New Roman"
$groupC.Font.Size = 12
$groupC.HorizontalAlignment = -4108
$groupC.VerticalAlignment = -4108
$groupC.WrapText = $true

# Group D: Times New Roman 12pt, centered/wrapped, no fill, no border
$groupD = $ws.Range("K4")
$groupD.Font.Name = "Times New Roman"
$groupD.Font.Size = 12
$groupD.HorizontalAlignment = -4108
$groupD.VerticalAlignment = -4108
$groupD.WrapText = $true

# Group F: Times New Roman 12pt, centered/wrapped, full medium box border
$groupF = $ws.Range("N2")
$groupF.Font.Name = "Times New Roman"
$groupF.Font.Size = 12
$groupF.HorizontalAlignment = -4108
$groupF.VerticalAlignment = -4108
$groupF.WrapText = $true
$groupF.Borders.Item(7).LineStyle = 1
$groupF.Borders.Item(7).Weight = -4138
$groupF.Borders.Item(7).Color = 0
$groupF.Borders.Item(8).LineStyle = 1
$groupF.Borders.Item(8).Weight = -4138
$groupF.Borders.Item(8).Color = 0
$groupF.Borders.Item(9).LineStyle = 1
$groupF.Borders.Item(9).Weight = -4138
$groupF.Borders.Item(9).Color = 0
$groupF.Borders.Item(10).LineStyle = 1
$groupF.Borders.Item(10).Weight = -4138
$groupF.Borders.Item(10).Color = 0

# Group G: Times New Roman 12pt, centered/wrapped, left/right/bottom medium border (no top)
$groupG = $ws.Range("N3,N4,N5,N6")
$groupG.Font.Name = "Times New Roman"
$groupG.Font.Size = 12
$groupG.HorizontalAlignment = -4108
$groupG.VerticalAlignment = -4108
$groupG.WrapText = $true
$groupG.Borders.Item(7).LineStyle = 1
$groupG.Borders.Item(7).Weight = -4138
$groupG.Borders.Item(7).Color = 0
$groupG.Borders.Item(9).LineStyle = 1
$groupG.Borders.Item(9).Weight = -4138
$groupG.Borders.Item(9).Color = 0
$groupG.Borders.Item(10).LineStyle = 1
$groupG.Borders.Item(10).Weight = -4138
$groupG.Borders.Item(10).Color = 0

# --- Selection -------------------------------------------------------------
$ws.Range("J6").Select()
